$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at 330; this shifts the existing rows 330-353
# down to 331-354 (matches the diff, which is a weekly-refresh row
# insertion at the top of this block, pushing the rest of the series down
# by one row).
$ws.Rows("330:330").Insert()

# Populate the newly inserted row 330 with the new weekly record.
$ws.Cells.Item(330, 1).Value = 8
$ws.Cells.Item(330, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(330, 3).Value = "Coquimbo"
$ws.Cells.Item(330, 4).Value = 44714
$ws.Cells.Item(330, 5).Value = 4
$ws.Cells.Item(330, 6).Value = 100114013
$ws.Cells.Item(330, 7).Value = "Zanahoria"
$ws.Cells.Item(330, 8).Value = "Sin especificar"
$ws.Cells.Item(330, 9).Value = "Primera"
$ws.Cells.Item(330, 10).Value = 560
$ws.Cells.Item(330, 11).Value = 6000
$ws.Cells.Item(330, 12).Value = 7000
$ws.Cells.Item(330, 13).Value = 6500
$ws.Cells.Item(330, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(330, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(330, 16).Value = 325
$ws.Cells.Item(330, 17).Value = 20
$ws.Cells.Item(330, 18).Value = "Hortaliza"
